$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "98÷6=16, 2"
$t.Cell(1, 2).Range.Text = "69÷7=9, 6"
$t.Cell(1, 3).Range.Text = "43÷2=21, 1"
$t.Cell(1, 4).Range.Text = "66÷8=8, 2"
$t.Cell(1, 5).Range.Text = "33÷7=4, 5"

$t.Cell(5, 1).Range.Text = "69÷9=7, 6"
$t.Cell(5, 2).Range.Text = "72÷9=8, 0"
$t.Cell(5, 3).Range.Text = "51÷4=12, 3"
$t.Cell(5, 4).Range.Text = "98÷4=24, 2"
$t.Cell(5, 5).Range.Text = "67÷3=22, 1"

$t.Cell(9, 1).Range.Text = "48÷5=9, 3"
$t.Cell(9, 2).Range.Text = "66÷7=9, 3"
$t.Cell(9, 3).Range.Text = "94÷3=31, 1"
$t.Cell(9, 4).Range.Text = "85÷5=17, 0"
$t.Cell(9, 5).Range.Text = "85÷3=28, 1"

$t.Cell(13, 1).Range.Text = "37÷5=7, 2"
$t.Cell(13, 2).Range.Text = "99÷7=14, 1"
$t.Cell(13, 3).Range.Text = "83÷3=27, 2"
$t.Cell(13, 4).Range.Text = "99÷6=16, 3"
$t.Cell(13, 5).Range.Text = "26÷5=5, 1"

$t.Cell(17, 1).Range.Text = "96÷7=13, 5"
$t.Cell(17, 2).Range.Text = "50÷6=8, 2"
$t.Cell(17, 3).Range.Text = "10÷4=2, 2"
$t.Cell(17, 4).Range.Text = "52÷6=8, 4"
$t.Cell(17, 5).Range.Text = "72÷2=36, 0"
